# Update the RACI matrix: column E ("I :" informed) values were revised
# (and one newly duplicated "C / I " label with a trailing space was introduced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value  = "R / A / C / I"
$ws.Range("E8").Value  = "R / A / C / I"
$ws.Range("E9").Value  = "R / A / C / I"
$ws.Range("E10").Value = "R / A / C / I"
$ws.Range("E11").Value = "C / I"
$ws.Range("E12").Value = "R / A / C / I"
$ws.Range("E13").Value = "R / A"
$ws.Range("E14").Value = "R / A"
$ws.Range("E15").Value = "C / I "
$ws.Range("E16").Value = "C / I"
$ws.Range("E17").Value = "C / I"
$ws.Range("E18").Value = "C / I"
$ws.Range("E19").Value = "C / I"
$ws.Range("E20").Value = "C / I"
$ws.Range("E21").Value = "C / I"
$ws.Range("E22").Value = "C / I"
$ws.Range("E23").Value = "C / I"
$ws.Range("E24").Value = "C / I"
$ws.Range("E25").Value = "C / I"
$ws.Range("E27").Value = "C / I"
$ws.Range("E28").Value = "C / I"
$ws.Range("E29").Value = "C / I"

# Move the active selection to E25, matching the saved view state.
$ws.Range("E25").Select()
